# Generate Report for Archive
# Update the status text "Ready for handoff" -> "In Translation" wherever it
# occurs (Overview sheet E2/F2, and the Status column on each language
# sheet), then autofit the affected columns since the new text is shorter
# than the old one.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $text = $cell.Text.ToString()
        if ($text -ceq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).AutoFit() | Out-Null

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).AutoFit() | Out-Null
